# Add a new "Delete.Name" test-case row to the "Location" worksheet,
# mirroring the existing "Edit.Name" row (row 16) / XPath-style rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location")

# Row 15 (".../Btn.Delete" row) carries the style (fill/border) we need to
# replicate onto the new row 17: s="22" for A:C, s="23" for D:G.
$ws.Range("A15:G15").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 20.25

# Populate the new row's values.
$ws.Range("A17").Value = "Delete.Name"
$ws.Range("B17").Value = "XPath"
$ws.Range("C17").Value = "//form/div/h4"
